# Update test-result cells on Sheet1 (rows 10-14): "Actual Outcome" (F) and
# "Fail/Pass" (G) columns, reflecting re-run test outputs / updated test cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 10: test now fails with a captured-stdin OSError instead of "-" ;
# Fail/Pass (G10) stays "Fail".
$ws.Range("F10").Value = "OSError: reading from stdin while output is captured"

# Rows 11-14: actual outcome now matches expected outcome, so these move
# from Fail to Pass.
$ws.Range("F11").Value = "Same as expected outcome."
$ws.Range("G11").Value = "Pass"

$ws.Range("F12").Value = "Same as expected outcome."
$ws.Range("G12").Value = "Pass"

$ws.Range("F13").Value = "Same as expected outcome."
$ws.Range("G13").Value = "Pass"

$ws.Range("F14").Value = "Same as expected outcome."
$ws.Range("G14").Value = "Pass"

# Update the view state to match where the author left the selection/scroll
# position on save.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("G11:G14").Select()
